# Test_Specifikation_Inloggningsfunktion.xlsx - add "Test Environment",
# "Version", "IDE", "Version " columns to the header row, in between the
# existing "Test case description" and "Expected result" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns at G (pushes the old G:I - "Expected result",
# "Actual result", "Verdict" - to K:M, and auto-expands the E2:M2 merged
# banner/title cell to E2:Q2).
$ws.Range("G1:J1").EntireColumn.Insert()

# Fill in the new header cells. Typed in this order (G, J, H, I) so the
# shared-string table is populated the same way the source workbook has
# it (Test Environment, Version<space>, Version, IDE).
$ws.Range("G4").Value = "Test Environment"
$ws.Range("J4").Value = "Version "
$ws.Range("H4").Value = "Version"
$ws.Range("I4").Value = "IDE"

# Re-set the narrower "Req link" column width and approximate the
# autosized widths of the new / shifted columns.
$ws.Columns("C").ColumnWidth = 7.451822916666667
$ws.Columns("G").ColumnWidth = 15.307291666666666
$ws.Columns("H").ColumnWidth = 6.877604166666667
$ws.Columns("I").ColumnWidth = 3.0221354166666665
$ws.Columns("J").ColumnWidth = 7.307291666666667
$ws.Columns("K").ColumnWidth = 13.877604166666666
$ws.Columns("L").ColumnWidth = 12.736979166666666
$ws.Columns("M").ColumnWidth = 11.307291666666666

# Match the author's final selection.
[void]$ws.Range("J5").Select()
